$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "142.9"
$ws.Range("C2").Value = "172.9 ⮟ -4"
$ws.Range("D2").Value = "195.9 ⮝ 4"
$ws.Range("I2").Value = 1.43
$ws.Range("J2").Value = "⮟"
$ws.Range("K2").Value = -4
$ws.Range("L2").Value = "⮟ -4"

$ws.Range("B3").Value = "157.9"
$ws.Range("C3").Value = "181.9 ⮟ -4"
$ws.Range("D3").Value = "219.9 ⮝ 4"
$ws.Range("I3").Value = 1.58
$ws.Range("J3").Value = "⮟"
$ws.Range("K3").Value = -4
$ws.Range("L3").Value = "⮟ -4"

$ws.Range("B4").Value = "158.9"
$ws.Range("D4").Value = "218.9 ⮝ 2"
$ws.Range("I4").Value = 1.59
$ws.Range("J4").Value = "⮟"
$ws.Range("K4").Value = -5
$ws.Range("L4").Value = "⮟ -5"

$ws.Range("B5").Value = "125.9"
$ws.Range("C5").Value = "150.9 ⮟ -2"
$ws.Range("D5").Value = "174.9 ⮝ 2"
$ws.Range("I5").Value = 1.26
$ws.Range("J5").Value = "⮟"
$ws.Range("K5").Value = -2
$ws.Range("L5").Value = "⮟ -2"

$ws.Range("B6").Value = "141.9"
$ws.Range("C6").Value = "171.9 ⮟ -4"
$ws.Range("D6").Value = "195.9 ⮝ 4"
$ws.Range("I6").Value = 1.42
$ws.Range("J6").Value = "⮟"
$ws.Range("K6").Value = -4
$ws.Range("L6").Value = "⮟ -4"

$ws.Range("B7").Value = "142.9"
$ws.Range("C7").Value = "172.9 ⮟ -4"
$ws.Range("D7").Value = "195.9 ⮝ 4"
$ws.Range("I7").Value = 1.43
$ws.Range("J7").Value = "⮟"
$ws.Range("K7").Value = -4
$ws.Range("L7").Value = "⮟ -4"

$ws.Range("B9").Value = "146.9"
$ws.Range("C9").Value = "172.9 ⮟ -2"
$ws.Range("I9").Value = 1.47
$ws.Range("J9").Value = "⮟"
$ws.Range("K9").Value = -2
$ws.Range("L9").Value = "⮟ -2"

$ws.Range("B10").Value = "123.9"
$ws.Range("C10").Value = "151.9 ⮟ -2"
$ws.Range("I10").Value = 1.24
$ws.Range("J10").Value = "⮟"
$ws.Range("K10").Value = -2
$ws.Range("L10").Value = "⮟ -2"

$ws.Range("B12").Value = "142.9"
$ws.Range("C12").Value = "172.9 ⮟ -4"
$ws.Range("D12").Value = "195.9 ⮝ 4"
$ws.Range("I12").Value = 1.43
$ws.Range("J12").Value = "⮟"
$ws.Range("K12").Value = -4
$ws.Range("L12").Value = "⮟ -4"

$ws.Range("B14").Value = "142.9"
$ws.Range("C14").Value = "172.9 ⮟ -4"
$ws.Range("D14").Value = "195.9 ⮝ 4"
$ws.Range("I14").Value = 1.43
$ws.Range("J14").Value = "⮟"
$ws.Range("K14").Value = -4
$ws.Range("L14").Value = "⮟ -4"

$ws.Range("B15").Value = "163.9"
$ws.Range("C15").Value = "185.9 ⮟ -2"
$ws.Range("D15").Value = "199.9 ⮟ -2"
$ws.Range("I15").Value = 1.64
$ws.Range("J15").Value = "⮟"
$ws.Range("K15").Value = -2
$ws.Range("L15").Value = "⮟ -2"

$ws.Range("B16").Value = "163.9"
$ws.Range("C16").Value = "185.9 ⮟ -2"
$ws.Range("D16").Value = "199.9 ⮟ -2"
$ws.Range("I16").Value = 1.64
$ws.Range("J16").Value = "⮟"
$ws.Range("K16").Value = -2
$ws.Range("L16").Value = "⮟ -2"

$ws.Range("B17").Value = "142.9"
$ws.Range("C17").Value = "166.9 ⮟ -2"
$ws.Range("D17").Value = "200.9 ⮝ 2"
$ws.Range("I17").Value = 1.43
$ws.Range("J17").Value = "⮟"
$ws.Range("K17").Value = -4
$ws.Range("L17").Value = "⮟ -4"

$ws.Range("B18").Value = "141.9"
$ws.Range("C18").Value = "171.9 ⮟ -4"
$ws.Range("D18").Value = "195.9 ⮝ 4"
$ws.Range("I18").Value = 1.42
$ws.Range("J18").Value = "⮟"
$ws.Range("K18").Value = -4
$ws.Range("L18").Value = "⮟ -4"

$ws.Range("B19").Value = "142.9"
$ws.Range("C19").Value = "172.9 ⮟ -4"
$ws.Range("D19").Value = "195.9 ⮝ 4"
$ws.Range("I19").Value = 1.43
$ws.Range("J19").Value = "⮟"
$ws.Range("K19").Value = -4
$ws.Range("L19").Value = "⮟ -4"

$ws.Range("B20").Value = "142.9"
$ws.Range("C20").Value = "172.9 ⮟ -4"
$ws.Range("D20").Value = "195.9 ⮝ 4"
$ws.Range("I20").Value = 1.43
$ws.Range("J20").Value = "⮟"
$ws.Range("K20").Value = -4
$ws.Range("L20").Value = "⮟ -4"

$ws.Range("B22").Value = "142.9"
$ws.Range("C22").Value = "172.9 ⮟ -4"
$ws.Range("D22").Value = "195.9 ⮝ 4"
$ws.Range("I22").Value = 1.43
$ws.Range("J22").Value = "⮟"
$ws.Range("K22").Value = -4
$ws.Range("L22").Value = "⮟ -4"

$ws.Range("B23").Value = "142.9"
$ws.Range("C23").Value = "172.9 ⮟ -4"
$ws.Range("D23").Value = "195.9 ⮝ 4"
$ws.Range("I23").Value = 1.43
$ws.Range("J23").Value = "⮟"
$ws.Range("K23").Value = -4
$ws.Range("L23").Value = "⮟ -4"

$ws.Range("B24").Value = "141.9"
$ws.Range("C24").Value = "172.9 ⮟ -3"
$ws.Range("D24").Value = "195.9 ⮝ 4"
$ws.Range("I24").Value = 1.42
$ws.Range("J24").Value = "⮟"
$ws.Range("K24").Value = -4
$ws.Range("L24").Value = "⮟ -4"

$ws.Range("B25").Value = "142.9"
$ws.Range("C25").Value = "172.9 ⮟ -21"
$ws.Range("D25").Value = "201.9 ⮟ -3"
$ws.Range("I25").Value = 1.43
$ws.Range("J25").Value = "⮟"
$ws.Range("K25").Value = -4
$ws.Range("L25").Value = "⮟ -4"

$ws.Range("B26").Value = "137.9"
$ws.Range("C26").Value = "162.9 ⮟ -2"
$ws.Range("D26").Value = "198.9 ⮝ 2"
$ws.Range("I26").Value = 1.38
$ws.Range("J26").Value = "⮟"
$ws.Range("K26").Value = -2
$ws.Range("L26").Value = "⮟ -2"

$ws.Range("B28").Value = "162.9"
$ws.Range("C28").Value = "184.9 ⮟ -3"
$ws.Range("I28").Value = 1.63
$ws.Range("J28").Value = "⮟"
$ws.Range("K28").Value = -4
$ws.Range("L28").Value = "⮟ -4"

$ws.Range("B29").Value = "143.9"
$ws.Range("C29").Value = "167.9 ⮟ -3"
$ws.Range("I29").Value = 1.44
$ws.Range("J29").Value = "⮟"
$ws.Range("K29").Value = -4
$ws.Range("L29").Value = "⮟ -4"

$ws.Range("B30").Value = "143.9"
$ws.Range("C30").Value = "167.9 ⮟ -3"
$ws.Range("I30").Value = 1.44
$ws.Range("J30").Value = "⮟"
$ws.Range("K30").Value = -4
$ws.Range("L30").Value = "⮟ -4"

$ws.Range("B31").Value = "142.9"
$ws.Range("C31").Value = "172.9 ⮟ -4"
$ws.Range("D31").Value = "195.9 ⮝ 4"
$ws.Range("I31").Value = 1.43
$ws.Range("J31").Value = "⮟"
$ws.Range("K31").Value = -4
$ws.Range("L31").Value = "⮟ -4"

$ws.Range("B34").Value = "151.9"
$ws.Range("C34").Value = "178.9 ⮟ -4"
$ws.Range("I34").Value = 1.52
$ws.Range("J34").Value = "⮟"
$ws.Range("K34").Value = -4
$ws.Range("L34").Value = "⮟ -4"

$ws.Range("B35").Value = "140.9"
$ws.Range("C35").Value = "167.9 ⮟ -3"
$ws.Range("I35").Value = 1.41
$ws.Range("J35").Value = "⮟"
$ws.Range("K35").Value = -3
$ws.Range("L35").Value = "⮟ -3"

$ws.Range("B36").Value = "160.9"
$ws.Range("C36").Value = "182.9 ⮟ -5"
$ws.Range("I36").Value = 1.61
$ws.Range("J36").Value = "⮟"
$ws.Range("K36").Value = -5
$ws.Range("L36").Value = "⮟ -5"

$ws.Range("B37").Value = "142.9"
$ws.Range("C37").Value = "172.9 ⮟ -4"
$ws.Range("D37").Value = "195.9 ⮝ 4"
$ws.Range("I37").Value = 1.43
$ws.Range("J37").Value = "⮟"
$ws.Range("K37").Value = -4
$ws.Range("L37").Value = "⮟ -4"

$ws.Range("B38").Value = "142.9"
$ws.Range("C38").Value = "172.9 ⮟ -4"
$ws.Range("D38").Value = "195.9 ⮝ 4"
$ws.Range("I38").Value = 1.43
$ws.Range("J38").Value = "⮟"
$ws.Range("K38").Value = -4
$ws.Range("L38").Value = "⮟ -4"

$ws.Range("B39").Value = "147.9"
$ws.Range("C39").Value = "165.9 ⮟ -4"
$ws.Range("I39").Value = 1.48
$ws.Range("J39").Value = "⮟"
$ws.Range("K39").Value = -4
$ws.Range("L39").Value = "⮟ -4"
